$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string "sdfd" -> "sdfdtyt" (used by cell E2)
$ws.Range("E2").Value = "sdfdtyt"

# Duplicate row 3 into row 4 (same values/types as row 3)
$ws.Range("A3:E3").Copy($ws.Range("A4"))
